# Fruta / hortaliza, semanal
#
# Insert a new weekly price record as row 8 in the daily logic sheet for
# "Agrícola del Norte S.A. de Arica - Durazno". Inserting the row pushes
# every existing record from the former row 8 down by one (through the
# former row 34, which becomes row 35), which matches how this workbook
# accumulates one new sample per update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (8..34) down to make room for the new sample.
$ws.Rows(8).Insert()

# Populate the newly opened row 8 with this week's sample.
$ws.Range("A8").Value  = 1
$ws.Range("B8").Value  = "Agrícola del Norte S.A. de Arica"
$ws.Range("C8").Value  = "Arica y Parinacota"
$ws.Range("D8").Value  = 44565
$ws.Range("E8").Value  = 15
$ws.Range("F8").Value  = "Fruta"
$ws.Range("G8").Value  = 100103
$ws.Range("H8").Value  = "Frutos de hueso (carozo)"
$ws.Range("I8").Value  = 100103004
$ws.Range("J8").Value  = "Durazno"
$ws.Range("K8").Value  = "Early Majestic"
$ws.Range("L8").Value  = "Segunda"
$ws.Range("M8").Value  = 270
$ws.Range("N8").Value  = 19000
$ws.Range("O8").Value  = 20000
$ws.Range("P8").Value  = 19500
$ws.Range("Q8").Value  = "$/bandeja 18 kilos granel"
$ws.Range("R8").Value  = "Región de O'Higgins"
$ws.Range("S8").Value  = 1083
$ws.Range("T8").Value  = 18
